$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.851.93"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "3.816.65"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'605.03"
$ws.Range("E5").Value = "  +1.53%  "

$ws.Range("D6").Value = "'166.21"
$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("E10").Value = "  +0.93%  "

$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("D12").Value = "'0.0000252"
$ws.Range("E12").Value = "  -0.72%  "

$ws.Range("D13").Value = "'36.13"
$ws.Range("E13").Value = "  +0.12%  "

$ws.Range("D14").Value = "4.451.37"
$ws.Range("E14").Value = "  +0.79%  "

$ws.Range("D15").Value = "3.817.50"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").Value = "67.856.25"
$ws.Range("E16").Value = "  +0.23%  "

$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("D18").Value = "'7.10"
$ws.Range("E18").Value = "  +1.09%  "

$ws.Range("E19").Value = "  +1.83%  "

$ws.Range("D20").Value = "'464.14"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("E21").Value = "  -1.33%  "

$ws.Range("D22").Value = "'0.703"
$ws.Range("E22").Value = "  +0.97%  "

$ws.Range("E23").Value = "  -4.38%  "

$ws.Range("D24").Value = "'83.37"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "'12.12"
$ws.Range("E25").Value = "  +1.34%  "

$ws.Range("D26").Value = "'2.13"
$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").Value = "'10.04"
$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("E28").Value = "  -0.14%  "

$ws.Range("D29").Value = "3.961.86"
$ws.Range("E29").Value = "  +0.87%  "

$ws.Range("D30").Value = "'2.80"
$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("E31").Value = "  +2.65%  "

$ws.Range("D32").Value = "'2.23"
$ws.Range("E32").Value = "  -1.06%  "

$ws.Range("D33").Value = "'29.62"
$ws.Range("E33").Value = "  -0.23%  "

$ws.Range("E34").Value = "  +0.19%  "

$ws.Range("D35").Value = "'9.10"
$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.57%  "

$ws.Range("D39").Value = "'5.83"
$ws.Range("E39").Value = "  +1.24%  "

$ws.Range("D40").Value = "'3.23"
$ws.Range("E40").Value = "  -4.07%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D43").Value = "'44.52"
$ws.Range("E43").Value = "  -2.82%  "

$ws.Range("D44").Value = "'47.73"
$ws.Range("E44").Value = "  -0.93%  "

$ws.Range("E45").Value = "  +0.60%  "

$ws.Range("D46").Value = "'28.04"
$ws.Range("E46").Value = "  +6.63%  "

$ws.Range("D47").Value = "'151.71"
$ws.Range("E47").Value = "  +1.51%  "

$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.38"
$ws.Range("E48").Value = "  +11.66%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'8.38"
$ws.Range("E49").Value = "  +0.66%  "

$ws.Range("D50").Value = "'1.86"
$ws.Range("E50").Value = "  +1.78%  "

$ws.Range("D51").Value = "'391.25"
$ws.Range("E51").Value = "  -0.67%  "
